$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: config-server
$ws.Range("A6").Value = "配置服务器"
$ws.Range("B6").Value = "config-server"
$ws.Range("C6").Value = 8888
$ws.Range("D6").Value = 8889
$ws.Rows.Item(6).RowHeight = 24.95
$ws.Range("C1").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)

# Row 7: config-client
$ws.Range("A7").Value = "配置客户端"
$ws.Range("B7").Value = "config-client"
$ws.Range("C7").Value = 8899
$ws.Range("D7").Value = 8900
$ws.Rows.Item(7).RowHeight = 24.95
$ws.Range("C1").Copy()
$ws.Range("C7:D7").PasteSpecial(-4122)

# Row 8: rabbitMQ (message queue) - highlighted in red instead of green
$ws.Range("A8").Value = "消息队列"
$ws.Range("B8").Value = "rabbitMQ"
$ws.Range("C8").Value = 5672
$ws.Rows.Item(8).RowHeight = 24.95
$ws.Range("C8").Interior.Color = 255

$ws.Range("D8").Select()
